$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.1295
$ws.Range("A8").Value = -21.16660000000001
$ws.Range("A10").Value = -20.45099999999998
$ws.Range("A12").Value = -22.27380000000003
$ws.Range("B13").Value = 6.571099999999997
$ws.Range("A18").Value = -22.25920000000002
$ws.Range("E20").Value = 12.0159

$wb.Save()
